# matplotlibWk5.1.pptx revision
#  1) Bump the "datetimeFigureOut" date placeholders from 12/7/18 -> 12/8/18
#     (handout master, notes master, slide master, and the two title-ish
#     slide layouts that carry their own date placeholder).
#  2) Fix a typo on slide 2 ("out plots" -> "our plots").

$p = $ppt.ActivePresentation

# ---- 1) Date placeholders ------------------------------------------------
$newDate = "12/8/18"

# Handout Master -> "Date Placeholder 2"
$hm = $p.HandoutMaster
$hm.Shapes.Item(2).TextFrame.TextRange.Text = $newDate

# Notes Master -> "Date Placeholder 2"
$nm = $p.NotesMaster
$nm.Shapes.Item(2).TextFrame.TextRange.Text = $newDate

# Slide Master -> "Date Placeholder 3"
$sm = $p.SlideMaster
$sm.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

# Slide Layouts (CustomLayouts) that carry their own date placeholder:
#   layout 6 -> "Date Placeholder 1"
#   layout 7 -> "Date Placeholder 3"
$layouts = $sm.CustomLayouts
$layouts.Item(6).Shapes.Item(1).TextFrame.TextRange.Text = $newDate
$layouts.Item(7).Shapes.Item(3).TextFrame.TextRange.Text = $newDate

# ---- 2) Slide 2 typo fix --------------------------------------------------
$s2 = $p.Slides.Item(2)
$body = $s2.Shapes.Item(2).TextFrame.TextRange
$para = $body.Paragraphs(5, 1)
$word = $para.Characters(80, 3)
$word.Text = "our"
